# Update countries & provincias Spain
# - Swap the Uruguay / Burkina Faso rows (alphabetical reorder)
# - Swap the Timor Oriental / Santa Lucia rows (alphabetical reorder)
# - Refresh the "last updated" timestamp
# - Refresh the per-country case numbers that changed between the two pulls

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- last-updated timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 29 de Septiembre de 2020 a las 01:36"

# --- reorder: Uruguay now comes before Burkina Faso ---
$ws.Range("A154").Value = "Uruguay"
$ws.Range("A155").Value = "Burkina Faso"

# --- reorder: Timor Oriental now comes before Santa Lucia ---
$ws.Range("A207").Value = "Timor Oriental"
$ws.Range("A208").Value = "Santa Lucia"

# --- refreshed numbers (Estados Unidos, row 4) ---
$ws.Range("B4").Value = 7357797
$ws.Range("C4").Value = 33604
$ws.Range("D4").Value = 4604514
$ws.Range("E4").Value = 2543541
$ws.Range("G4").Value = 289
$ws.Range("H4").Value = 209742

# --- refreshed numbers (Brasil, row 6) ---
$ws.Range("B6").Value = 4748327
$ws.Range("C6").Value = 16018
$ws.Range("D6").Value = 4084182
$ws.Range("E6").Value = 521984
$ws.Range("G6").Value = 385
$ws.Range("H6").Value = 142161

# --- refreshed numbers (Colombia, row 8) ---
$ws.Range("B8").Value = 818203
$ws.Range("C8").Value = 5147
$ws.Range("D8").Value = 722536
$ws.Range("E8").Value = 70026
$ws.Range("G8").Value = 153
$ws.Range("H8").Value = 25641

# --- refreshed numbers (Peru, row 9) ---
$ws.Range("B9").Value = 808714
$ws.Range("C9").Value = 3412
$ws.Range("D9").Value = 670989
$ws.Range("E9").Value = 105401
$ws.Range("G9").Value = 62
$ws.Range("H9").Value = 32324

# --- refreshed numbers (row 25) ---
$ws.Range("B25").Value = 288618
$ws.Range("C25").Value = 2280
$ws.Range("D25").Value = 252400
$ws.Range("E25").Value = 26673

# --- refreshed numbers (row 29) ---
$ws.Range("B29").Value = 155301
$ws.Range("C29").Value = 2176
$ws.Range("D29").Value = 132607
$ws.Range("E29").Value = 13416
$ws.Range("G29").Value = 10
$ws.Range("H29").Value = 9278

# --- refreshed numbers (row 38) ---
$ws.Range("B38").Value = 111277
$ws.Range("C38").Value = 722
$ws.Range("D38").Value = 87695
$ws.Range("E38").Value = 21234
$ws.Range("G38").Value = 8
$ws.Range("H38").Value = 2348

# --- refreshed numbers (row 57) ---
$ws.Range("B57").Value = 65883
$ws.Range("C57").Value = 1286
$ws.Range("D57").Value = 31799
$ws.Range("E57").Value = 33466
$ws.Range("G57").Value = 12
$ws.Range("H57").Value = 618

# --- refreshed numbers (row 58) ---
$ws.Range("B58").Value = 58460
$ws.Range("C58").Value = 136
$ws.Range("D58").Value = 49895
$ws.Range("E58").Value = 7454
$ws.Range("G58").Value = 3
$ws.Range("H58").Value = 1111

# --- refreshed numbers (row 93) ---
$ws.Range("B93").Value = 13788
$ws.Range("C93").Value = 90
$ws.Range("E93").Value = 2324

# --- refreshed numbers (row 114) ---
$ws.Range("B114").Value = 7474
$ws.Range("C114").Value = 10
$ws.Range("D114").Value = 7091
$ws.Range("E114").Value = 222

# --- refreshed numbers (row 130) ---
$ws.Range("B130").Value = 4797
$ws.Range("C130").Value = 79
$ws.Range("D130").Value = 1813
$ws.Range("E130").Value = 2808
$ws.Range("G130").Value = 2
$ws.Range("H130").Value = 176

# --- refreshed numbers (row 133) ---
$ws.Range("B133").Value = 4386
$ws.Range("C133").Value = 24
$ws.Range("E133").Value = 2039

# --- refreshed numbers (row 154, now Uruguay) ---
$ws.Range("B154").Value = 2010
$ws.Range("C154").Value = 2
$ws.Range("D154").Value = 1755
$ws.Range("E154").Value = 208
$ws.Range("H154").Value = 47

# --- refreshed numbers (row 155, now Burkina Faso) ---
$ws.Range("D155").Value = 1276
$ws.Range("E155").Value = 676
$ws.Range("H155").Value = 56

# --- refreshed numbers (row 169) ---
$ws.Range("D169").Value = 885
$ws.Range("E169").Value = 11

# --- refreshed numbers (row 185) ---
$ws.Range("D185").Value = 315
$ws.Range("E185").Value = 1

# --- refreshed numbers (row 190) ---
$ws.Range("B190").Value = 211
$ws.Range("C190").Value = 1
$ws.Range("D190").Value = 208
